$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "'676"
$ws.Range("F2").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F8").Value = 1
